$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Highlight wat we want..." -> "Highlight what we want..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Highlight wat we", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Highlight what we", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. The paragraph that used to read "Click graph white space to see options"
#    becomes "Flip data  if needed " and a new sub-bullet paragraph
#    ("Click graph white space to see options") is added right after it.
# ---------------------------------------------------------------------------
$pClick = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd("`r`a") -eq "Click graph white space to see options") {
        $pClick = $d.Paragraphs.Item($i)
        break
    }
}

$clickRange = $pClick.Range
$clickRange.Collapse(0)
$clickRange.InsertParagraphAfter() | Out-Null

$pClickIndex = $pClick.Range.ListFormat.ListLevelNumber
$newClickPara = $d.Paragraphs.Item($pClick.Range.Information(1) )

# Re-locate paragraphs by scanning again (indices shifted after insertion)
$pFlip = $null
$pNewClick = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -eq "Click graph white space to see options`r") {
        $pFlip = $d.Paragraphs.Item($i)
        $pNewClick = $d.Paragraphs.Item($i + 1)
        break
    }
}

$pNewClick.Range.Text = "Click graph white space to see options"
$pNewClick.Range.ListFormat.ListLevelNumber = $pFlip.Range.ListFormat.ListLevelNumber + 1
$pFlip.Range.Text = "Flip data  if needed "

# ---------------------------------------------------------------------------
# 3. After "Select 1 under horizontal", append the rest of the graphing
#    instructions as new list paragraphs.
# ---------------------------------------------------------------------------
$pSelect = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Select 1 under horizontal`r") {
        $pSelect = $d.Paragraphs.Item($i)
        break
    }
}
$selectIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Select 1 under horizontal`r") {
        $selectIndex = $i
        break
    }
}

$newItems = @(
  @{t="Add labels"; lvl=3},
  @{t="Click green " + [char]8220 + "+" + [char]8221; lvl=4},
  @{t="Turn on axis titles"; lvl=4},
  @{t="Rename graph and axes titles"; lvl=4},
  @{t="Manually fix scaling issues"; lvl=3},
  @{t="Double click white space on graph to bring up graph- options sidebar"; lvl=4},
  @{t="Click on either set of axis labels to make " + [char]8220 + "axis options" + [char]8221 + " graphic appear in sidebar"; lvl=4},
  @{t="Under " + [char]8220 + "Axis Options>Axis Options" + [char]8221 + " change the bounds to appropriate limits"; lvl=4},
  @{t="Iff this is the first time you are reaching step 4 this for this graph repeat steps 1-3 for the other axis"; lvl=4},
  @{t="Add trend line and equation"; lvl=3},
  @{t="Click green " + [char]8220 + "+" + [char]8221; lvl=4},
  @{t="Check the box next to " + [char]8220 + "Trendline" + [char]8221; lvl=4},
  @{t="Click arrow next to " + [char]8220 + "Trendline" + [char]8221; lvl=4},
  @{t="Scroll to near bottom and check the appropriate box"; lvl=4},
  @{t="Make adjustments as required and appropriate"; lvl=4}
)

$anchor = $pSelect.Range
$anchor.Collapse(0)
for ($k = 0; $k -lt $newItems.Count; $k++) {
    $anchor.InsertParagraphAfter() | Out-Null
}

for ($k = 0; $k -lt $newItems.Count; $k++) {
    $para = $d.Paragraphs.Item($selectIndex + 1 + $k)
    $para.Range.Text = $newItems[$k].t
    $para.Range.ListFormat.ListLevelNumber = $newItems[$k].lvl
}

Write-Host "Final paragraph count:" $d.Paragraphs.Count
